# Revert "Predicting PM10 hourly levels / Organizing Data Gathering and Exploration code"
#
# This undoes the earlier commit that:
#   - renamed sheet "Sheet1" -> "arima_graph"
#   - added a new "pm10_limits" worksheet with PM10 limit reference data
#
# so we revert by:
#   - renaming "arima_graph" back to "Sheet1"
#   - repointing the chart series formulas from the old sheet name to the new one
#   - removing the "pm10_limits" worksheet that was added

$wb = $excel.ActiveWorkbook

# Rename the remaining worksheet back to its original name
$ws1 = $wb.Worksheets.Item("arima_graph")
$ws1.Name = "Sheet1"

# Update the chart's series formulas so they reference the renamed sheet
$co = $ws1.ChartObjects(1)
$chart = $co.Chart
$s1 = $chart.SeriesCollection(1)
$s1.Formula = "=SERIES(Sheet1!`$B`$1,,Sheet1!`$B`$2:`$B`$13,1)"
$s2 = $chart.SeriesCollection(2)
$s2.Formula = "=SERIES(Sheet1!`$C`$1,,Sheet1!`$C`$2:`$C`$13,2)"

# Remove the "pm10_limits" sheet that was added by the commit being reverted
$excel.DisplayAlerts = $false
$ws2 = $wb.Worksheets.Item("pm10_limits")
$ws2.Delete()
$excel.DisplayAlerts = $true
